$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at O:P (pushing former O..U to Q..W)
$ws.Columns("O:P").Insert()

# Rename the (now shifted-in-place) headers for M and N
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# New column headers for the inserted O and P columns
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# New column data values (rows 2-6) for the inserted O and P columns
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1

$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 4

$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 2

$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1

$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 2
